$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 currently holds "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)\n"
# Row 24 currently holds "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)\n"
# The edit swaps these two entries so row 23 now shows the LOT2038 text
# and row 24 shows the LOT2028 text.

$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$lot2038 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"

$ws.Range("B23").Value = $lot2038
$ws.Range("C23").Value = $lot2038

$ws.Range("B24").Value = $lot2028
$ws.Range("C24").Value = $lot2028
